# Applies the "added nl-core-careplan profile and mapping + adjusted mappings
# for NursingIntervention and OutcomeOfCare" edit to the Data sheet (and
# the selection on that sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Update the "Mapping to FHIR" / "Equivalence" / "Notes" columns ---

# Row 3 (OutcomeOfCare rootconcept mapping)
$ws.Range("P3").Value = "CarePlan  / Procedure"
$ws.Range("R3").Value = "** OutcomeOfCare`nMaps to CarePlan / or Procedure and referenced resources`nAdd mappings to `n* CarePlan`n* NursingIntervention`n* Procedure`n* TextResult`n* GeneralMeasurement`n* FunctionalOrMentalStatus`n"

# Row 4 (HealthcareResult mapping)
$ws.Range("P4").Value = "DiagnosticReport.conclusion (HCIM Textresult)"

# Row 7 (NursingIntervention mapping)
$ws.Range("P7").Value = "Procedure / Careplan.activity"
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()

# --- Adjust column widths ---
$ws.Range("O1").ColumnWidth = 13.85546875
$ws.Range("R1").ColumnWidth = 34.28515625

# --- Adjust row heights ---
$ws.Range("B3").RowHeight = 178.5
$ws.Range("B6").RowHeight = 89.25
$ws.Range("B7").RowHeight = 89.25

# --- Move the active selection from P7 to P6 ---
$ws.Range("P6").Select()
